# Add toggle to control key limiter in CLS and main menu.
$wb = $excel.ActiveWorkbook

# --- KeyLimiter sheet: add two new rows for the new toggles ---
$keyLimiter = $wb.Worksheets.Item("KeyLimiter")

$keyLimiter.Cells.Item(8, 1).Value = "LIMIT_CLS"
$keyLimiter.Cells.Item(8, 2).Value = "Limit keys in CLS (Custom Level Select)"

$keyLimiter.Cells.Item(9, 1).Value = "LIMIT_MAIN_MENU"
$keyLimiter.Cells.Item(9, 2).Value = "Limit keys in main menu"

# --- KeyViewer sheet: reword the Korean description string ---
$keyViewer = $wb.Worksheets.Item("KeyViewer")
$keyViewer.Cells.Item(3, 3).Value = "등록된 키들의 키뷰어를 보여줍니다."
